# Refresh cryptos list (prices & 1h change %) to match latest coinranking scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.583.28'
$ws.Range('E2').Value = '  +2.19%  '
$ws.Range('D3').Value = '1.663.51'
$ws.Range('D4').Value = "'0.9997"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'236.17"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D7').Value = "'0.4793"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').Value = "'0.2618"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = "'0.06155"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.41%  '
$ws.Range('E10').Value = '  +0.18%  '
$ws.Range('D11').Value = '1.666.76'
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('D12').Value = "'14.74"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.50%  '
$ws.Range('D13').Value = "'0.5906"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.00%  '
$ws.Range('D14').Value = "'4.378"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.98%  '
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').Value = '25.582.44'
$ws.Range('E18').Value = '  +2.23%  '
$ws.Range('D19').Value = "'0.000006752"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.17%  '
$ws.Range('E20').Value = '  +0.62%  '
$ws.Range('D21').Value = '1.875.68'
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('D22').Value = "'4.427"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').Value = "'8.648"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.46%  '
$ws.Range('D24').Value = "'5.307"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.53%  '
$ws.Range('D25').Value = "'134.55"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = "'15.06"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.55%  '
$ws.Range('D27').Value = "'1.405"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.98%  '
$ws.Range('D28').Value = "'104.63"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.14%  '
$ws.Range('D29').Value = "'1.686"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('E30').Value = '  +4.71%  '
$ws.Range('D31').Value = "'3.654"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.72%  '
$ws.Range('E32').Value = '  -3.17%  '
$ws.Range('D33').Value = "'0.9997"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').Value = "'0.04319"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.94%  '
$ws.Range('D35').Value = "'2.619"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.50%  '
$ws.Range('D36').Value = "'0.6117"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.16%  '
$ws.Range('D37').Value = "'0.9499"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.48%  '
$ws.Range('D38').Value = "'2.610"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('D39').Value = "'0.8540"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.35%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('B41').Value = 'PaxosStandard'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range('D41').Value = "'1.000"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('E42').Value = '  -1.85%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = "'1.878"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.40%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = "'97.90"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = "'0.3764"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.11%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = "'4.703"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.29%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = "'0.1119"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = "'6.210"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.50%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = "'0.05263"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.24%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = "'29.48"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('E51').Value = '  +0.09%  '
